# 10Th - MB for single stock and added new group
#
# MarketBeat weekly rating-history sheet: a new reporting week (Jun_27) has
# rolled in, so three new date columns are inserted in front of the existing
# ones (shifting Jun_17/Jun_15/Jun_13/Jun_10 three slots to the right), and
# two brand-new analyst/benchmark rows (Benchmark, Evercore ISI) are appended
# at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert three new date columns before the existing "B" column ---------
# This shifts old B,C,D,E (Jun_17,Jun_15,Jun_13,Jun_10) to E,F,G,H and keeps
# their custom column width tagging along for the ride.
$ws.Columns("B:D").Insert()

# New header cells: the newest week first (Jun_26 is written into both the
# second and third new columns, matching the source data), Jun_27 last so it
# lands at the end of the shared-string table.
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- Fill the three new rating columns for every existing analyst row -----
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# --- Add the two new rows for the newly tracked firms ----------------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

# Write the newest week's header last, after everything else, so it is the
# final entry appended to the shared-string table.
$ws.Range("B1").Value = "Jun_27"

# --- Re-assert the custom width (8 characters) across the whole date block -
$w = 7.166666666666667
$ws.Columns("C").ColumnWidth = $w
$ws.Columns("D").ColumnWidth = $w
$ws.Columns("E").ColumnWidth = $w
$ws.Columns("F").ColumnWidth = $w
$ws.Columns("G").ColumnWidth = $w
$ws.Columns("H").ColumnWidth = $w
